$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the answer-position strings: drop the stray "0," so
# "(x,y,0,0.02)" becomes "(x,y,0.02)". The correct-answer cells
# (already "(0,0.375,0.02)") and the question-text cells in column A
# are left untouched. Column order below (B, then C, then D, then E)
# matches the order these distinct values were (re)introduced into the
# shared-string table in the target workbook.

$ws.Range("B3").Value = "(0.25,0.25,0.02)"
$ws.Range("B4").Value = "(0.25,0.25,0.02)"
$ws.Range("B5").Value = "(0.25,0.25,0.02)"

$ws.Range("C2").Value = "(-0.25,0.25,0.02)"
$ws.Range("C4").Value = "(-0.25,0.25,0.02)"
$ws.Range("C5").Value = "(-0.25,0.25,0.02)"

$ws.Range("D2").Value = "(0.25,0.5,0.02)"
$ws.Range("D3").Value = "(0.25,0.5,0.02)"
$ws.Range("D4").Value = "(0.25,0.5,0.02)"

$ws.Range("E2").Value = "(-0.25,0.5,0.02)"
$ws.Range("E3").Value = "(-0.25,0.5,0.02)"
$ws.Range("E5").Value = "(-0.25,0.5,0.02)"

# --- Move the active selection from C17 to C7
$ws.Range("C7").Select() | Out-Null
